$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. All source cells are stored as
# text (inline strings) in the workbook, so every value is written as text
# (NumberFormat "@" keeps numeric-looking strings from being reinterpreted as
# numbers, and resetting the Style back to Normal avoids leaving a stray
# "Text" number-format style behind).
$updates = @{
    'D2' = '237.34'
    'D3' = '21.65'
    'B4' = 'HuobiToken'
    'C4' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D4' = '5.358'
    'E4' = '3HuobiTokenHT'
    'B5' = 'Cronos'
    'C5' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D5' = '0.05558'
    'E5' = '4CronosCRO'
    'B6' = 'GateToken'
    'C6' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D6' = '3.368'
    'E6' = '5GateTokenGT'
    'B7' = 'KuCoinToken'
    'C7' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'D7' = '6.454'
    'E7' = '6KuCoinTokenKCS'
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D8' = '0.8009'
    'E8' = '7MXTokenMX'
    'B9' = 'FTXToken'
    'C9' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D9' = '1.035'
    'E9' = '8FTXTokenFTT'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D10' = '0.1398'
    'E10' = '9WazirXWRX'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '0.07312'
    'E11' = '10MandalaExchangeTokenMDX'
    'B12' = 'LiechtensteinCryptoassetsExchange'
    'C12' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D12' = '0.03289'
    'E12' = '11LiechtensteinCryptoassetsExchangeLCX'
    'B13' = 'ProBitToken'
    'C13' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'D13' = '0.1287'
    'E13' = '12ProBitTokenPROB'
    'D14' = '0.02871'
    'D15' = '0.09245'
    'D16' = '0.001659'
    'D17' = '3.252'
    'D18' = '0.04754'
    'D19' = '0.0005709'
    'E19' = '18OneONE'
    'D20' = '0.006259'
    'D22' = '0.001052'
    'D23' = '0.0001500'
    'D24' = '0.0004181'
    'B25' = 'LEO'
    'C25' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D25' = '3.957'
    'E25' = '24LEOLEOBestin24h'
    'B26' = 'BTSEToken'
    'C26' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D26' = '2.200'
    'E26' = '25BTSETokenBTSE'
    'B27' = 'BitpandaEcosystemToken'
    'C27' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'D27' = '0.3289'
    'E27' = '26BitpandaEcosystemTokenBEST'
    'D40' = '0.04134'
    'D41' = '0.007012'
    'D42' = '0.002909'
    'D43' = '0.1037'
    'D44' = '0.008783'
    'D45' = '0.00005440'
    'D47' = '0.6799'
    'D48' = '0.03164'
    'D49' = '0.00002100'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}

Write-Host "Updated $($updates.Count) cells"
